# The edit inserts a new record row into the price table (a weekly price
# report gained a new entry). This shifts every existing row from 78
# downward down by one, and the row that used to be row 77 is duplicated
# into the new row 78, while row 77 itself is updated with the newly
# reported record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 78; rows 78..129 shift down to 79..130.
$ws.Rows(78).Insert()

# The row that is now at 77 (old data, untouched by the insert) needs to be
# duplicated down into the newly-blanked row 78.
$ws.Range("A77:T77").Copy()
$ws.Range("A78").PasteSpecial()

# Now overwrite row 77 with the new record's values.
$ws.Cells.Item(77, 4).Value2 = 45001
$ws.Cells.Item(77, 13).Value2 = 80
$ws.Cells.Item(77, 14).Value2 = 4000
$ws.Cells.Item(77, 15).Value2 = 4000
$ws.Cells.Item(77, 16).Value2 = 4000
$ws.Cells.Item(77, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(77, 19).Value2 = 2000
